$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column (D) sometimes holds values that look like plain decimal
# numbers (e.g. "231.52"); Excel would otherwise silently convert them to
# numeric cells. Force text formatting first so they stay text, matching the
# original inlineStr/string cell type.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '28.698.39'
$ws.Range("E2").Value = '  -1.70%  '
$ws.Range("D3").Value = '1.802.47'
$ws.Range("E3").Value = '  -1.32%  '
$ws.Range("E4").Value = '  +0.40%  '
$ws.Range("D5").Value = '231.52'
$ws.Range("E5").Value = '  -1.98%  '
$ws.Range("D6").Value = '0.5945'
$ws.Range("E6").Value = '  -2.95%  '
$ws.Range("D8").Value = '0.2776'
$ws.Range("E8").Value = '  -1.44%  '
$ws.Range("D9").Value = '0.06828'
$ws.Range("E9").Value = '  -3.93%  '
$ws.Range("D10").Value = '23.34'
$ws.Range("E10").Value = '  -1.06%  '
$ws.Range("D11").Value = '0.07544'
$ws.Range("E11").Value = '  -1.58%  '
$ws.Range("D12").Value = '1.809.11'
$ws.Range("E12").Value = '  -0.84%  '
$ws.Range("D13").Value = '4.785'
$ws.Range("E13").Value = '  -0.72%  '
$ws.Range("D14").Value = '0.6241'
$ws.Range("E14").Value = '  -1.42%  '
$ws.Range("D15").Value = '2.047.59'
$ws.Range("E15").Value = '  -1.27%  '
$ws.Range("D16").Value = '0.000009312'
$ws.Range("E16").Value = '  -7.66%  '
$ws.Range("D17").Value = '75.34'
$ws.Range("E17").Value = '  -4.60%  '
$ws.Range("D18").Value = '28.680.27'
$ws.Range("E18").Value = '  -1.67%  '
$ws.Range("D19").Value = '5.478'
$ws.Range("E19").Value = '  -6.71%  '
$ws.Range("D20").Value = '1.004'
$ws.Range("E20").Value = '  +0.36%  '
$ws.Range("D21").Value = '210.53'
$ws.Range("E21").Value = '  -7.57%  '
$ws.Range("E22").Value = '  -2.86%  '
$ws.Range("D23").Value = '6.854'
$ws.Range("E23").Value = '  -2.25%  '
$ws.Range("E24").Value = '  +0.38%  '
$ws.Range("D25").Value = '154.32'
$ws.Range("E25").Value = '  -0.52%  '
$ws.Range("D26").Value = '7.850'
$ws.Range("E26").Value = '  -2.40%  '
$ws.Range("D27").Value = '0.1274'
$ws.Range("E27").Value = '  -3.14%  '
$ws.Range("D28").Value = '16.39'
$ws.Range("E28").Value = '  -1.30%  '
$ws.Range("D29").Value = '1.431'
$ws.Range("E29").Value = '  -3.89%  '
$ws.Range("D30").Value = '0.06163'
$ws.Range("E30").Value = '  -3.21%  '
$ws.Range("D31").Value = '1.420'
$ws.Range("D32").Value = '3.783'
$ws.Range("E32").Value = '  -1.10%  '
$ws.Range("D33").Value = '3.748'
$ws.Range("E33").Value = '  -1.46%  '
$ws.Range("D34").Value = '1.719'
$ws.Range("E34").Value = '  -1.62%  '
$ws.Range("E35").Value = '  -6.01%  '
$ws.Range("D36").Value = '0.6401'
$ws.Range("E36").Value = '  -1.52%  '
$ws.Range("D37").Value = '2.495'
$ws.Range("E37").Value = '  -1.99%  '
$ws.Range("D38").Value = '2.713'
$ws.Range("E38").Value = '  -1.28%  '
$ws.Range("D39").Value = '0.01713'
$ws.Range("E39").Value = '  -1.70%  '
$ws.Range("D40").Value = '6.436'
$ws.Range("E40").Value = '  -2.20%  '
$ws.Range("D41").Value = '1.132.64'
$ws.Range("E41").Value = '  -6.99%  '
$ws.Range("D42").Value = '0.8687'
$ws.Range("E42").Value = '  -5.46%  '
$ws.Range("D43").Value = '1.007'
$ws.Range("E43").Value = '  +0.73%  '
$ws.Range("D44").Value = '100.66'
$ws.Range("E44").Value = '  -0.55%  '
$ws.Range("D45").Value = '1.966.21'
$ws.Range("E45").Value = '  -0.41%  '
$ws.Range("D46").Value = '60.57'
$ws.Range("E46").Value = '  -3.92%  '
$ws.Range("D47").Value = '0.00000000113'
$ws.Range("E47").Value = '  -5.07%  '
$ws.Range("D48").Value = '1.600'
$ws.Range("E48").Value = '  -1.52%  '
$ws.Range("D49").Value = '8.337'
$ws.Range("E49").Value = '  -2.97%  '
$ws.Range("D50").Value = '0.05472'
$ws.Range("E50").Value = '  -0.85%  '
$ws.Range("D51").Value = '0.4495'
$ws.Range("E51").Value = '  -1.46%  '

# Restore the default (Normal) style so the edited cells keep the same
# style index as the rest of the sheet (no explicit style attribute).
$ws.Range("D2:D51").Style = "Normal"
